$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.679.39'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '2.556.19'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.63'
$ws.Range("E5").Value = '  -2.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.75'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("E7").Value = '  -1.09%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.56'
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0807'
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.39'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '2.950.22'
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.87'
$ws.Range("E15").Value = '  +5.04%  '
$ws.Range("D16").Value = '2.578.33'
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.835'
$ws.Range("E17").Value = '  -1.75%  '
$ws.Range("D18").Value = '42.726.04'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.72'
$ws.Range("E19").Value = '  -1.46%  '
$ws.Range("D20").Value = '0.0₃0955'
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.31'
$ws.Range("E21").Value = '  -3.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.47'
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '247.22'
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.72'
$ws.Range("E26").Value = '  +1.03%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.94'
$ws.Range("E29").Value = '  -2.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.14'
$ws.Range("E30").Value = '  -2.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.06'
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.73'
$ws.Range("E32").Value = '  -3.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0794'
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("E35").Value = '  -3.63%  '
$ws.Range("E36").Value = '  -3.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.58'
$ws.Range("E37").Value = '  -1.50%  '
$ws.Range("E38").Value = '  +11.07%  '
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.117'
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.61'
$ws.Range("E41").Value = '  +1.88%  '
$ws.Range("E42").Value = '  +5.99%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("D45").Value = '1.987.32'
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.18'
$ws.Range("E46").Value = '  -2.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.95'
$ws.Range("E47").Value = '  -1.64%  '
$ws.Range("D48").Value = '2.803.39'
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '81.04'
$ws.Range("E49").Value = '  -4.02%  '
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.13'
$ws.Range("E51").Value = '  -2.72%  '
